$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.98254127752276
$ws.Range("C2").Value = 9.850571624199638
$ws.Range("E2").Value = 12.19072082085392
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.631538631868254
$ws.Range("I2").Value = 22.01393554855276
$ws.Range("L2").Value = 10.02503709274763
$ws.Range("M2").Value = 14.35175561457324
$ws.Range("N2").Value = 17.86030925221476
$ws.Range("O2").Value = 21.77956766534952

$ws.Range("B3").Value = 13.54443789954035
$ws.Range("C3").Value = 9.65440839314342
$ws.Range("E3").Value = 12.22798948911169
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.633417501913303
$ws.Range("I3").Value = 22.12079238304468
$ws.Range("L3").Value = 10.03276545845033
$ws.Range("M3").Value = 14.26592772848859
$ws.Range("N3").Value = 17.90679865464352
$ws.Range("O3").Value = 21.85062232810751

$ws.Range("B4").Value = 13.26978232397225
$ws.Range("C4").Value = 9.530766744907659
$ws.Range("E4").Value = 12.25218772206988
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.634632882231919
$ws.Range("I4").Value = 22.19097896259568
$ws.Range("L4").Value = 10.03890817888998
$ws.Range("M4").Value = 14.21497488706972
$ws.Range("N4").Value = 17.93715574901138
$ws.Range("O4").Value = 21.89964870356809

$ws.Range("B5").Value = 13.15660507975945
$ws.Range("C5").Value = 9.479619075653263
$ws.Range("E5").Value = 12.26238014655682
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.635143734427384
$ws.Range("I5").Value = 22.2207303309112
$ws.Range("L5").Value = 10.0417633183369
$ws.Range("M5").Value = 14.19466631171338
$ws.Range("N5").Value = 17.94998311483081
$ws.Range("O5").Value = 21.92098084310432

$ws.Range("B6").Value = 13.1377416380488
$ws.Range("C6").Value = 9.471081220379739
$ws.Range("E6").Value = 12.26409263166539
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.635229503146026
$ws.Range("I6").Value = 22.22573994015199
$ws.Range("L6").Value = 10.04225868177832
$ws.Range("M6").Value = 14.19132203497025
$ws.Range("N6").Value = 17.95214069232513
$ws.Range("O6").Value = 21.92460465554918

$ws.Range("B7").Value = 13.26826081468773
$ws.Range("C7").Value = 9.530079981860498
$ws.Range("E7").Value = 12.25232383764402
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.63463970863678
$ws.Range("I7").Value = 22.19137554577798
$ws.Range("L7").Value = 10.03894525870098
$ws.Range("M7").Value = 14.2146991349372
$ws.Range("N7").Value = 17.93732689343958
$ws.Range("O7").Value = 21.89993092146603

$ws.Range("B8").Value = 13.83275551104719
$ws.Range("C8").Value = 9.783622040489554
$ws.Range("E8").Value = 12.20329855676654
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.632173679822512
$ws.Range("I8").Value = 22.04982946705534
$ws.Range("L8").Value = 10.02741212231032
$ws.Range("M8").Value = 14.32180966824151
$ws.Range("N8").Value = 17.87596314532059
$ws.Range("O8").Value = 21.802944492379

$ws.Range("B9").Value = 14.88803830808115
$ws.Range("C9").Value = 10.25375427240892
$ws.Range("E9").Value = 12.11756078185608
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.62782553665321
$ws.Range("I9").Value = 21.80860841866458
$ws.Range("L9").Value = 10.01585895465579
$ws.Range("M9").Value = 14.54497074970486
$ws.Range("N9").Value = 17.76997133928864
$ws.Range("O9").Value = 21.65576134982656

$ws.Range("B10").Value = 15.62369069832963
$ws.Range("C10").Value = 10.58054290847376
$ws.Range("E10").Value = 12.06086185793745
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.624925211696802
$ws.Range("I10").Value = 21.65360101590003
$ws.Range("L10").Value = 10.01407871058138
$ws.Range("M10").Value = 14.71590591531043
$ws.Range("N10").Value = 17.70078999062576
$ws.Range("O10").Value = 21.57406833397363

$ws.Range("B11").Value = 15.94828492514819
$ws.Range("C11").Value = 10.72476100525097
$ws.Range("E11").Value = 12.03642431180536
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.623669019036333
$ws.Range("I11").Value = 21.58792196459065
$ws.Range("L11").Value = 10.0147160522342
$ws.Range("M11").Value = 14.79495841089427
$ws.Range("N11").Value = 17.67119384385566
$ws.Range("O11").Value = 21.54268940383227

$ws.Range("B12").Value = 16.06964907397194
$ws.Range("C12").Value = 10.7787034027041
$ws.Range("E12").Value = 12.02736455681903
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.623202367174917
$ws.Range("I12").Value = 21.56374756406458
$ws.Range("L12").Value = 10.01516455393859
$ws.Range("M12").Value = 14.82505997193758
$ws.Range("N12").Value = 17.66025535200775
$ws.Range("O12").Value = 21.53164189224698

$ws.Range("B13").Value = 16.04358175682463
$ws.Range("C13").Value = 10.76711621142754
$ws.Range("E13").Value = 12.02930710964843
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.623302467598631
$ws.Range("I13").Value = 21.56892293559522
$ws.Range("L13").Value = 10.01505876384702
$ws.Range("M13").Value = 14.81856999356541
$ws.Range("N13").Value = 17.66259920177747
$ws.Range("O13").Value = 21.53398398415498

$ws.Range("B14").Value = 15.95830131166606
$ws.Range("C14").Value = 10.72921247595586
$ws.Range("E14").Value = 12.03567507164959
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.623630446361948
$ws.Range("I14").Value = 21.58591914709495
$ws.Range("L14").Value = 10.01474880479718
$ws.Range("M14").Value = 14.79743167016978
$ws.Range("N14").Value = 17.67028854161151
$ws.Range("O14").Value = 21.54176376547868

$ws.Range("B15").Value = 15.90585942013227
$ws.Range("C15").Value = 10.70590721252173
$ws.Range("E15").Value = 12.0396009055577
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.623832518930951
$ws.Range("I15").Value = 21.59642061442187
$ws.Range("L15").Value = 10.01458589391311
$ws.Range("M15").Value = 14.78450487407803
$ws.Range("N15").Value = 17.67503348699706
$ws.Range("O15").Value = 21.5466379431458

$ws.Range("B16").Value = 15.60226595485087
$ws.Range("C16").Value = 10.57102583926247
$ws.Range("E16").Value = 12.06248612133669
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.625008574398471
$ws.Range("I16").Value = 21.65799071814152
$ws.Range("L16").Value = 10.01406610065742
$ws.Range("M16").Value = 14.71076397869897
$ws.Range("N16").Value = 17.70276183668544
$ws.Range("O16").Value = 21.57623569293862

$ws.Range("B17").Value = 15.41336573799128
$ws.Range("C17").Value = 10.48712028364903
$ws.Range("E17").Value = 12.07687206622165
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.625746196660614
$ws.Range("I17").Value = 21.69700137444332
$ws.Range("L17").Value = 10.01411723451915
$ws.Range("M17").Value = 14.66584349407851
$ws.Range("N17").Value = 17.72025198097839
$ws.Range("O17").Value = 21.5958767329536

$ws.Range("B18").Value = 15.3037742311453
$ws.Range("C18").Value = 10.43844387795491
$ws.Range("E18").Value = 12.08527406569697
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.626176406864488
$ws.Range("I18").Value = 21.71989422644118
$ws.Range("L18").Value = 10.01428292328083
$ws.Range("M18").Value = 14.6401293503551
$ws.Range("N18").Value = 17.73048835268635
$ws.Range("O18").Value = 21.60771793977738

$ws.Range("B19").Value = 15.26651022956033
$ws.Range("C19").Value = 10.42189240441898
$ws.Range("E19").Value = 12.08814077091071
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.626323091812088
$ws.Range("I19").Value = 21.7277234431188
$ws.Range("L19").Value = 10.01436245603623
$ws.Range("M19").Value = 14.63144469821923
$ws.Range("N19").Value = 17.73398454750336
$ws.Range("O19").Value = 21.61182053010901

$ws.Range("B20").Value = 15.43357269631086
$ws.Range("C20").Value = 10.49609547810674
$ws.Range("E20").Value = 12.07532745802257
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.625667060093865
$ws.Range("I20").Value = 21.69280152583516
$ws.Range("L20").Value = 10.01409769358554
$ws.Range("M20").Value = 14.67061277666703
$ws.Range("N20").Value = 17.7183718627127
$ws.Range("O20").Value = 21.59372956307567

$ws.Range("B21").Value = 15.98339317361105
$ws.Range("C21").Value = 10.74036414124368
$ws.Range("E21").Value = 12.03379938105297
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.623533866027643
$ws.Range("I21").Value = 21.58090802279224
$ws.Range("L21").Value = 10.01483423316107
$ws.Range("M21").Value = 14.80363615972606
$ws.Range("N21").Value = 17.66802270274002
$ws.Range("O21").Value = 21.53945596802882

$ws.Range("B22").Value = 16.3336413994495
$ws.Range("C22").Value = 10.89609033838594
$ws.Range("E22").Value = 12.00779011298101
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.622192378248093
$ws.Range("I22").Value = 21.51184160534117
$ws.Range("L22").Value = 10.01652255755479
$ws.Range("M22").Value = 14.89153341989194
$ws.Range("N22").Value = 17.63668382997196
$ws.Range("O22").Value = 21.50885345651115

$ws.Range("B23").Value = 16.14757133607752
$ws.Range("C23").Value = 10.81334459471037
$ws.Range("E23").Value = 12.02156840633536
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.622903549982071
$ws.Range("I23").Value = 21.548331386929
$ws.Range("L23").Value = 10.01551136624391
$ws.Range("M23").Value = 14.8445398833629
$ws.Range("N23").Value = 17.65326679148329
$ws.Range("O23").Value = 21.52474010188327

$ws.Range("B24").Value = 15.42444020666517
$ws.Range("C24").Value = 10.49203915454476
$ws.Range("E24").Value = 12.07602536627348
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.625702818602256
$ws.Range("I24").Value = 21.69469882872176
$ws.Range("L24").Value = 10.01410610344007
$ws.Range("M24").Value = 14.66845623598441
$ws.Range("N24").Value = 17.71922130014834
$ws.Range("O24").Value = 21.5946985877498

$ws.Range("B25").Value = 14.60898351234407
$ws.Range("C25").Value = 10.12969575777674
$ws.Range("E25").Value = 12.13964664334913
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.628949929110067
$ws.Range("I25").Value = 21.8699678135618
$ws.Range("L25").Value = 10.017803454597
$ws.Range("M25").Value = 14.48329866559744
$ws.Range("N25").Value = 17.79711499821067
$ws.Range("O25").Value = 21.6909500827273
